$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The new "Pollinaria removed" row belongs right after the
# "Inflorescences" row, i.e. right before the "Danaus plexippus
# abundance" row. Locate that row by its label so the insertion point
# doesn't depend on a hard-coded row index.
$beforeRow = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    if ($t.Rows.Item($i).Cells.Item(1).Range.Text -like "Danaus plexippus abundance*") {
        $beforeRow = $t.Rows.Item($i)
        break
    }
}

$newRow = $t.Rows.Add($beforeRow)

# Match the row height used in the target revision (w:trHeight val=572,
# hRule=auto). Row.Height is in points, OOXML w:val is in twentieths of
# a point, so 572/20 = 28.6.
$newRow.HeightRule = 0
$newRow.Height = 28.6

$newRow.Cells.Item(1).Range.Text = "Pollinaria removed"
$newRow.Cells.Item(2).Range.Text = "1"
$newRow.Cells.Item(3).Range.Text = "0.156"
$newRow.Cells.Item(4).Range.Text = "0.406"
$newRow.Cells.Item(5).Range.Text = "0.172"
$newRow.Cells.Item(6).Range.Text = "0.376"
$newRow.Cells.Item(7).Range.Text = "0.154"
$newRow.Cells.Item(8).Range.Text = "0.406"
$newRow.Cells.Item(9).Range.Text = "0.171"
$newRow.Cells.Item(10).Range.Text = "0.374"
